$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Range("A1").EntireColumn.Insert()

# Populate the new column A
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "19th"

# Rename the sheet
$ws.Name = "Navdeep Saini"
